# "process refactoring the output"
# Rename the worksheet "Uncut_Sheet" -> "Uncut_Sheet_1" and keep the
# Print_Area defined name (and its sheet-qualified reference) in sync,
# then update the active selection on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the sheet.
$ws.Name = "Uncut_Sheet_1"

# 2. Re-apply the print area so the _xlnm.Print_Area defined name is
#    rewritten to reference the new sheet name (Uncut_Sheet_1!$A$1:$G$42).
$ws.PageSetup.PrintArea = '$A$1:$G$42'

# 3. Update the current selection on the sheet.
[void]$ws.Range("B14:C14").Select()
